$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 18 (current REFREG99 row), shifting it down to row 19
$ws.Rows.Item(18).Insert()

# Fill the new row 18 with the "Niet te lokaliseren" entry
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "REFREG93"
$ws.Range("C18").Value = "Niet te lokaliseren"
$ws.Range("D18").Value = "Niet te lokaliseren"

# Update the volgnr of the row that got shifted down (now row 19)
$ws.Range("A19").Value = 18
